$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 24205
$ws.Range("E2").Value = -760
$ws.Range("F2").Value = 459
$ws.Range("G2").Value = -572
$ws.Range("H2").Value = 423
$ws.Range("I2").Value = -199
$ws.Range("J2").Value = 622
$ws.Range("K2").Value = 74187
$ws.Range("L2").Value = 41632
$ws.Range("M2").Value = 32555
$ws.Range("N2").Value = 27867
$ws.Range("O2").Value = 4688
$ws.Range("P2").Value = 1272
$ws.Range("Q2").Value = 1711
$ws.Range("R2").Value = -3431
$ws.Range("S2").Value = 3361
$ws.Range("T2").Value = 6978
$ws.Range("U2").Value = -5267
$ws.Range("V2").Value = 27383
$ws.Range("W2").Value = -3.14
$ws.Range("X2").Value = 1.75
$ws.Range("Y2").Value = -0.71
$ws.Range("Z2").Value = 0.57
$ws.Range("AA2").Value = 127.88
$ws.Range("AB2").Value = 2101.11
$ws.Range("AC2").Value = -835
$ws.Range("AD2").Value = -94.13
$ws.Range("AE2").Value = 116845
$ws.Range("AF2").Value = 0.67
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 0.25
$ws.Range("AI2").Value = -23.95
$ws.Range("AJ2").Value = 23849371
$ws.Range("D3").Value = 23015
$ws.Range("E3").Value = -1446
$ws.Range("F3").Value = -1446
$ws.Range("G3").Value = -3031
$ws.Range("H3").Value = 1821
$ws.Range("I3").Value = 1003
$ws.Range("J3").Value = 818
$ws.Range("K3").Value = 72988
$ws.Range("L3").Value = 40566
$ws.Range("M3").Value = 32422
$ws.Range("N3").Value = 29553
$ws.Range("O3").Value = 2869
$ws.Range("P3").Value = 1272
$ws.Range("Q3").Value = 205
$ws.Range("R3").Value = -1971
$ws.Range("S3").Value = 1683
$ws.Range("T3").Value = 8760
$ws.Range("U3").Value = -8555
$ws.Range("V3").Value = 26275
$ws.Range("W3").Value = -6.28
$ws.Range("X3").Value = 7.91
$ws.Range("Y3").Value = 3.49
$ws.Range("Z3").Value = 2.47
$ws.Range("AA3").Value = 125.12
$ws.Range("AB3").Value = 2208.55
$ws.Range("AC3").Value = 4206
$ws.Range("AD3").Value = 17.83
$ws.Range("AE3").Value = 123914
$ws.Range("AF3").Value = 0.61
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 23849371
$ws.Range("D4").Value = 27367
$ws.Range("E4").Value = 1325
$ws.Range("F4").Value = 1214
$ws.Range("G4").Value = -1035
$ws.Range("H4").Value = 2194
$ws.Range("I4").Value = 2421
$ws.Range("J4").Value = -227
$ws.Range("K4").Value = 62486
$ws.Range("L4").Value = 29840
$ws.Range("M4").Value = 32646
$ws.Range("N4").Value = 31957
$ws.Range("O4").Value = 689
$ws.Range("P4").Value = 1272
$ws.Range("Q4").Value = 4076
$ws.Range("R4").Value = -1393
$ws.Range("S4").Value = -4039
$ws.Range("T4").Value = 4418
$ws.Range("U4").Value = -342
$ws.Range("V4").Value = 22442
$ws.Range("W4").Value = 4.84
$ws.Range("X4").Value = 8.02
$ws.Range("Y4").Value = 7.87
$ws.Range("Z4").Value = 3.24
$ws.Range("AA4").Value = 91.40000000000001
$ws.Range("AB4").Value = 2391.3
$ws.Range("AC4").Value = 10152
$ws.Range("AD4").Value = 7.75
$ws.Range("AE4").Value = 133995
$ws.Range("AF4").Value = 0.59
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 0.51
$ws.Range("AI4").Value = 3.94
$ws.Range("AJ4").Value = 23849371
$ws.Range("D5").Value = 36316
$ws.Range("E5").Value = 2844
$ws.Range("F5").Value = 2844
$ws.Range("G5").Value = 2797
$ws.Range("H5").Value = 2326
$ws.Range("I5").Value = 2349
$ws.Range("J5").Value = -23
$ws.Range("K5").Value = 60778
$ws.Range("L5").Value = 26621
$ws.Range("M5").Value = 34157
$ws.Range("N5").Value = 33523
$ws.Range("O5").Value = 634
$ws.Range("P5").Value = 1272
$ws.Range("Q5").Value = 4803
$ws.Range("R5").Value = 4157
$ws.Range("S5").Value = -2951
$ws.Range("T5").Value = 1127
$ws.Range("U5").Value = 3676
$ws.Range("V5").Value = 19161
$ws.Range("W5").Value = 7.83
$ws.Range("X5").Value = 6.41
$ws.Range("Y5").Value = 7.17
$ws.Range("Z5").Value = 3.77
$ws.Range("AA5").Value = 77.94
$ws.Range("AB5").Value = 2564.28
$ws.Range("AC5").Value = 9849
$ws.Range("AD5").Value = 13.81
$ws.Range("AE5").Value = 140563
$ws.Range("AF5").Value = 0.97
$ws.Range("AG5").Value = 1950
$ws.Range("AH5").Value = 1.43
$ws.Range("AI5").Value = 19.8
$ws.Range("AJ5").Value = 23849371
$ws.Range("D6").Value = 31121
$ws.Range("E6").Value = 1587
$ws.Range("F6").Value = 1587
$ws.Range("G6").Value = 1305
$ws.Range("H6").Value = 1038
$ws.Range("I6").Value = 1038
$ws.Range("K6").Value = 56596
$ws.Range("L6").Value = 21670
$ws.Range("M6").Value = 34926
$ws.Range("N6").Value = 34294
$ws.Range("P6").Value = 1272
$ws.Range("Q6").Value = 6765
$ws.Range("R6").Value = -3989
$ws.Range("S6").Value = -4602
$ws.Range("T6").Value = 2767
$ws.Range("U6").Value = 3999
$ws.Range("V6").Value = 15088
$ws.Range("W6").Value = 5.1
$ws.Range("X6").Value = 3.34
$ws.Range("Y6").Value = 3.06
$ws.Range("Z6").Value = 1.77
$ws.Range("AA6").Value = 62.05
$ws.Range("AB6").Value = 2604.63
$ws.Range("AC6").Value = 4352
$ws.Range("AD6").Value = 24.59
$ws.Range("AE6").Value = 143796
$ws.Range("AF6").Value = 0.74
$ws.Range("AG6").Value = 850
$ws.Range("AH6").Value = 0.79
$ws.Range("AI6").Value = 19.53
$ws.Range("AJ6").Value = 23849371
$ws.Range("D7").Value = 26195
$ws.Range("E7").Value = -1578
$ws.Range("G7").Value = -2128
$ws.Range("H7").Value = -1846
$ws.Range("I7").Value = -1794
$ws.Range("K7").Value = 53537
$ws.Range("L7").Value = 20543
$ws.Range("M7").Value = 32993
$ws.Range("N7").Value = 32407
$ws.Range("P7").Value = 1270
$ws.Range("Q7").Value = 1900
$ws.Range("R7").Value = -2540
$ws.Range("S7").Value = -1033
$ws.Range("T7").Value = 2540
$ws.Range("U7").Value = -1015
$ws.Range("W7").Value = -6.02
$ws.Range("X7").Value = -7.05
$ws.Range("Y7").Value = -5.38
$ws.Range("Z7").Value = -3.35
$ws.Range("AA7").Value = 62.27
$ws.Range("AC7").Value = -7522
$ws.Range("AD7").Value = -7.42
$ws.Range("AE7").Value = 135881
$ws.Range("AF7").Value = 0.41
$ws.Range("AG7").Value = 133
$ws.Range("AH7").Value = 0.24
$ws.Range("AI7").Value = -1.77
$ws.Range("D8").Value = 28622
$ws.Range("E8").Value = -187
$ws.Range("G8").Value = -69
$ws.Range("H8").Value = -66
$ws.Range("I8").Value = -69
$ws.Range("K8").Value = 53247
$ws.Range("L8").Value = 20507
$ws.Range("M8").Value = 32740
$ws.Range("N8").Value = 32150
$ws.Range("P8").Value = 1270
$ws.Range("Q8").Value = 3287
$ws.Range("R8").Value = -2603
$ws.Range("S8").Value = -643
$ws.Range("T8").Value = 3000
$ws.Range("U8").Value = 1510
$ws.Range("W8").Value = -0.65
$ws.Range("X8").Value = -0.23
$ws.Range("Y8").Value = -0.21
$ws.Range("Z8").Value = -0.12
$ws.Range("AA8").Value = 62.63
$ws.Range("AC8").Value = -288
$ws.Range("AD8").Value = -193.81
$ws.Range("AE8").Value = 134804
$ws.Range("AF8").Value = 0.41
$ws.Range("AG8").Value = 367
$ws.Range("AH8").Value = 0.66
$ws.Range("AI8").Value = -127.35
$ws.Range("D9").Value = 29163
$ws.Range("E9").Value = 819
$ws.Range("G9").Value = 722
$ws.Range("H9").Value = 578
$ws.Range("I9").Value = 567
$ws.Range("K9").Value = 53323
$ws.Range("L9").Value = 20263
$ws.Range("M9").Value = 33060
$ws.Range("N9").Value = 32460
$ws.Range("P9").Value = 1270
$ws.Range("Q9").Value = 3713
$ws.Range("R9").Value = -2043
$ws.Range("S9").Value = -753
$ws.Range("T9").Value = 2500
$ws.Range("U9").Value = 2135
$ws.Range("W9").Value = 2.81
$ws.Range("X9").Value = 1.98
$ws.Range("Y9").Value = 1.76
$ws.Range("Z9").Value = 1.08
$ws.Range("AA9").Value = 61.29
$ws.Range("AC9").Value = 2379
$ws.Range("AD9").Value = 23.46
$ws.Range("AE9").Value = 136104
$ws.Range("AF9").Value = 0.41
$ws.Range("AG9").Value = 500
$ws.Range("AH9").Value = 0.9
$ws.Range("AI9").Value = 21.02
